$d = $word.ActiveDocument

# Locate the empty, bold (sz 28) paragraph that sits right after the
# "Git commit" explanation and right before the "Git push" section.
# It is uniquely identifiable as the empty paragraph whose mark is
# bold with a 28 half-point (14pt) size.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    $rng = $cand.Range
    if ($rng.Text.Trim() -eq "" -and $rng.Font.Bold -eq -1 -and $rng.Font.Size -eq 14) {
        $target = $cand
        break
    }
}

$prev = $target.Previous()

# Drop the old (bold / 14pt) paragraph mark entirely, then grow a brand
# new paragraph off the preceding (non-bold, 12pt) paragraph so the new
# paragraph mark naturally carries the surrounding, non-bold formatting
# instead of inheriting the bold/size-28 run properties we are removing.
$target.Range.Delete()
$prev.Range.InsertParagraphAfter()

$newPara = $prev.Next()
$newRange = $newPara.Range
$newRange.Text = "We can also see the commit history online by clicking the commit buttons with a clock icon"

$newRange = $newPara.Range
$newRange.Font.Size = 12
$newRange.Font.SizeBi = 12
$newRange.LanguageID = "cs-CZ"
